# This script regenerates the "K" column (column G, formerly computed from a
# "Strike#" quantity) plus the two downstream cells (H12, I12) that were
# recalculated alongside it, as described in the commit message:
#   "regen save_data to use K instead of Strike#, regen std/mean,
#    calc and write s_vals"
#
# The workbook already contains the old values; we overwrite the cells with
# the freshly (re)computed values, exactly as the regenerated save_data
# pipeline would when writing the .xlsx back out.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$updates = @{
    "G2"  = 1
    "G3"  = 1
    "G4"  = 1
    "G5"  = 3
    "G6"  = 0
    "G7"  = 1
    "G8"  = 2
    "G9"  = 2
    "G10" = 0
    "G11" = 1
    "G12" = 3
    "H12" = 2
    "I12" = 7
    "G13" = 0
    "G14" = 0
    "G15" = 0
    "G17" = 1
    "G18" = 2
    "G19" = 1
    "G20" = 2
    "G21" = 1
    "G22" = 0
    "G23" = 0
    "G24" = 0
    "G25" = 1
    "G26" = 0
    "G27" = 1
    "G28" = 2
    "G29" = 1
    "G30" = 1
    "G31" = 0
    "G32" = 0
    "G33" = 0
    "G34" = 0
    "G35" = 2
    "G36" = 0
    "G37" = 3
    "G38" = 1
    "G39" = 0
    "G40" = 0
    "G41" = 0
    "G42" = 1
    "G43" = 1
    "G44" = 4
    "G45" = 3
    "G46" = 2
    "G47" = 2
    "G48" = 1
    "G49" = 0
    "G50" = 0
    "G51" = 1
    "G52" = 1
    "G54" = 0
    "G55" = 2
    "G56" = 0
    "G57" = 0
    "G58" = 0
    "G59" = 1
    "G60" = 3
    "G61" = 1
    "G62" = 3
    "G63" = 5
    "G64" = 2
    "G65" = 0
    "G66" = 0
    "G67" = 2
    "G68" = 1
    "G69" = 0
    "G70" = 1
    "G71" = 3
    "G72" = 0
    "G73" = 0
    "G74" = 0
    "G75" = 2
    "G76" = 4
    "G78" = 2
    "G79" = 1
    "G80" = 1
    "G81" = 1
    "G82" = 2
    "G83" = 1
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
